$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue 'D2' '29.126.71'
Set-TextValue 'E2' '  -0.22%  '
Set-TextValue 'D3' '1.842.07'
Set-TextValue 'E3' '  -0.37%  '
Set-TextValue 'D4' '0.9994'
Set-TextValue 'E4' '  -0.06%  '
Set-TextValue 'D5' '241.16'
Set-TextValue 'E5' '  -2.06%  '
Set-TextValue 'D6' '0.6870'
Set-TextValue 'E6' '  -1.56%  '
Set-TextValue 'D7' '1.0000'
Set-TextValue 'E7' '  -0.03%  '
Set-TextValue 'D8' '0.3021'
Set-TextValue 'E8' '  -1.28%  '
Set-TextValue 'D9' '0.07466'
Set-TextValue 'E9' '  -3.26%  '
Set-TextValue 'D10' '23.15'
Set-TextValue 'E10' '  -1.49%  '
Set-TextValue 'D11' '0.07670'
Set-TextValue 'E11' '  -1.95%  '
Set-TextValue 'B12' 'WrappedEther'
Set-TextValue 'C12' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D12' '1.843.77'
Set-TextValue 'E12' '  -0.21%  '
Set-TextValue 'B13' 'Polkadot'
Set-TextValue 'C13' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D13' '5.067'
Set-TextValue 'E13' '  -1.24%  '
Set-TextValue 'D14' '0.6833'
Set-TextValue 'E14' '  -0.39%  '
Set-TextValue 'D15' '87.47'
Set-TextValue 'E15' '  -5.98%  '
Set-TextValue 'D16' '6.176'
Set-TextValue 'E16' '  -7.06%  '
Set-TextValue 'D17' '29.107.12'
Set-TextValue 'E17' '  -0.29%  '
Set-TextValue 'D18' '0.000008163'
Set-TextValue 'E18' '  -1.85%  '
Set-TextValue 'D19' '2.079.72'
Set-TextValue 'E19' '  -0.58%  '
Set-TextValue 'D20' '228.31'
Set-TextValue 'E20' '  -5.39%  '
Set-TextValue 'D21' '12.55'
Set-TextValue 'E21' '  -1.47%  '
Set-TextValue 'D22' '0.9997'
Set-TextValue 'E22' '  +0.01%  '
Set-TextValue 'D23' '7.397'
Set-TextValue 'E23' '  -1.66%  '
Set-TextValue 'E24' '  -0.04%  '
Set-TextValue 'D25' '0.1457'
Set-TextValue 'E25' '  -3.46%  '
Set-TextValue 'D26' '159.99'
Set-TextValue 'E26' '  +0.65%  '
Set-TextValue 'E27' '  -0.82%  '
Set-TextValue 'E28' '  -1.07%  '
Set-TextValue 'D29' '1.513'
Set-TextValue 'E29' '  -2.18%  '
Set-TextValue 'D30' '4.265'
Set-TextValue 'E30' '  +0.82%  '
Set-TextValue 'D31' '4.138'
Set-TextValue 'E31' '  -0.84%  '
Set-TextValue 'D32' '1.198'
Set-TextValue 'E32' '  +0.50%  '
Set-TextValue 'D33' '0.05200'
Set-TextValue 'E33' '  +1.56%  '
Set-TextValue 'D34' '0.7660'
Set-TextValue 'E34' '  -3.91%  '
Set-TextValue 'E35' '  -1.21%  '
Set-TextValue 'E36' '  -1.12%  '
Set-TextValue 'D37' '2.679'
Set-TextValue 'E37' '  -0.49%  '
Set-TextValue 'D38' '1.316.16'
Set-TextValue 'E38' '  +0.02%  '
Set-TextValue 'D39' '0.01837'
Set-TextValue 'E39' '  -1.93%  '
Set-TextValue 'D40' '2.728'
Set-TextValue 'E40' '  +0.54%  '
Set-TextValue 'D41' '0.9323'
Set-TextValue 'E41' '  -1.16%  '
Set-TextValue 'D42' '104.79'
Set-TextValue 'E42' '  -2.19%  '
Set-TextValue 'D43' '5.774'
Set-TextValue 'E43' '  -4.05%  '
Set-TextValue 'D44' '0.9995'
Set-TextValue 'E44' '  -0.05%  '
Set-TextValue 'D45' '1.980.93'
Set-TextValue 'E45' '  -0.39%  '
Set-TextValue 'D46' '0.5196'
Set-TextValue 'E46' '  +0.31%  '
Set-TextValue 'E47' '  -0.32%  '
Set-TextValue 'D48' '64.95'
Set-TextValue 'E48' '  +1.33%  '
Set-TextValue 'D49' '9.504'
Set-TextValue 'E49' '  -2.40%  '
Set-TextValue 'D50' '1.774'
Set-TextValue 'E50' '  +0.59%  '
Set-TextValue 'B51' 'XinFinNetwork'
Set-TextValue 'C51' 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
Set-TextValue 'D51' '0.07379'
Set-TextValue 'E51' '  +16.72%  '
